# "convert logic to import books by isbn"
#
# The original sheet ("Names") had two columns - Name / Email - with the
# email column rendered as a mailto: hyperlink. This converts the sheet
# into an "Import" sheet with a single ISBN column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Names" -> "Import"
$ws.Name = "Import"

# Drop the mailto: hyperlinks that lived on column B before clearing the
# column itself - Hyperlinks.Delete() removes the whole collection (and
# the now-orphaned relationships) in one shot.
$ws.Hyperlinks.Delete()

# Column B (Email) goes away entirely - clear its values/number formats/
# styles so the sheet is back down to a single A column.
$ws.Range("B:B").Clear()

# The "Hyperlink" cell style (used by the old email column) is no longer
# referenced by anything, so remove it from the workbook's style gallery.
$wb.Styles.Item("Hyperlink").Delete()

# Header + two ISBNs (write A3 before A2 so the shared-string table comes
# out in the same order as the authored workbook).
$ws.Range("A1").Value = "ISBN"
$ws.Range("A3").Value = "978-1338216660"
$ws.Range("A2").Value = "978-0590353427"

# Matches the saved cursor position in the authored workbook.
$ws.Range("B34").Select()
